$d = $word.ActiveDocument

# Locate the last paragraph in the document (the one containing "Vikram")
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)

# Insert a new paragraph after it and set its text to "karthick".
# This mirrors the run/paragraph formatting (en-US language) of the
# preceding paragraph, matching the target diff.
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($lastParaIndex + 1)
$newPara.Range.Text = "karthick"
